# Relayout - round 5 - Zip and Unzip
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phieu Danh Gia")

# Fill column E (TỰ ĐÁNH GIÁ) with "ok" for rows 9 through 20
$ws.Range("E9:E20").Value = "ok"

# Update the active selection to reflect where the user ended up editing
$ws.Activate()
$ws.Range("E16").Select()
